$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44533
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 2000
$ws.Range("L2").Value = 2200
$ws.Range("M2").Value = 2100
$ws.Range("P2").Value = 2100

# Row 3
$ws.Range("D3").Value = 44532
$ws.Range("J3").Value = 60
$ws.Range("K3").Value = 2000
$ws.Range("L3").Value = 2200
$ws.Range("M3").Value = 2100
$ws.Range("P3").Value = 2100

# Row 4
$ws.Range("D4").Value = 44610
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 600
$ws.Range("L4").Value = 650
$ws.Range("M4").Value = 625
$ws.Range("P4").Value = 625

# Row 5
$ws.Range("D5").Value = 44624
$ws.Range("J5").Value = 120
$ws.Range("K5").Value = 650
$ws.Range("L5").Value = 700
$ws.Range("M5").Value = 675
$ws.Range("P5").Value = 675

# Row 6
$ws.Range("D6").Value = 44608
$ws.Range("J6").Value = 120
$ws.Range("K6").Value = 600
$ws.Range("L6").Value = 650
$ws.Range("M6").Value = 625
$ws.Range("P6").Value = 625
